$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 3
$ws.Range("G2").Value = 126.9318136666667
$ws.Range("H2").Value = 380.795441
$ws.Range("I2").Value = 0.1973293860115714
$ws.Range("J2").Value = 0.1973293860115715
$ws.Range("K2").Value = 3
$ws.Range("M2").Value = 24.455837
$ws.Range("N2").Value = 73.36751100000001
$ws.Range("O2").Value = 0.1553502885444182
$ws.Range("P2").Value = 0.1553502885444182
$ws.Range("Q2").Value = 3104.223745146372
$ws.Range("R2").Value = 27938.01370631735
$ws.Range("S2").Value = 0.0306551770551905
$ws.Range("T2").Value = 0.03065517705519051

$ws.Range("E3").Value = 3
$ws.Range("G3").Value = 126.9318136666667
$ws.Range("H3").Value = 380.795441
$ws.Range("I3").Value = 0.1973293860115714
$ws.Range("J3").Value = 0.1973293860115715
$ws.Range("K3").Value = 3
$ws.Range("M3").Value = 33.36516466666667
$ws.Range("N3").Value = 100.095494
$ws.Range("O3").Value = 0.2119448194841458
$ws.Range("P3").Value = 0.2119448194841458
$ws.Range("Q3").Value = 4235.100864426983
$ws.Range("R3").Value = 38115.90777984286
$ws.Range("S3").Value = 0.04182294109713983
$ws.Range("T3").Value = 0.04182294109713984

$ws.Range("E4").Value = 3
$ws.Range("G4").Value = 126.9318136666667
$ws.Range("H4").Value = 380.795441
$ws.Range("I4").Value = 0.1973293860115714
$ws.Range("J4").Value = 0.1973293860115715
$ws.Range("K4").Value = 3
$ws.Range("M4").Value = 43.331795
$ws.Range("N4").Value = 129.995385
$ws.Range("O4").Value = 0.2752556314632608
$ws.Range("P4").Value = 0.2752556314632608
$ws.Range("Q4").Value = 5500.183328782198
$ws.Range("R4").Value = 49501.64995903978
$ws.Range("S4").Value = 0.05431602475287264
$ws.Range("T4").Value = 0.05431602475287264

$ws.Range("E5").Value = 3
$ws.Range("G5").Value = 126.9318136666667
$ws.Range("H5").Value = 380.795441
$ws.Range("I5").Value = 0.1973293860115714
$ws.Range("J5").Value = 0.1973293860115715
$ws.Range("K5").Value = 3
$ws.Range("M5").Value = 45.91844699999999
$ws.Range("N5").Value = 137.755341
$ws.Range("O5").Value = 0.2916867654524183
$ws.Range("P5").Value = 0.2916867654524183
$ws.Range("Q5").Value = 5828.511758466708
$ws.Range("R5").Value = 52456.60582620037
$ws.Range("S5").Value = 0.05755837033442695
$ws.Range("T5").Value = 0.05755837033442696

$ws.Range("E6").Value = 3
$ws.Range("G6").Value = 126.9318136666667
$ws.Range("H6").Value = 380.795441
$ws.Range("I6").Value = 0.1973293860115714
$ws.Range("J6").Value = 0.1973293860115715
$ws.Range("K6").Value = 3
$ws.Range("M6").Value = 10.35258366666667
$ws.Range("N6").Value = 31.057751
$ws.Range("O6").Value = 0.06576249505575693
$ws.Range("P6").Value = 0.06576249505575693
$ws.Range("Q6").Value = 1314.07222094591
$ws.Range("R6").Value = 11826.64998851319
$ws.Range("S6").Value = 0.01297687277194152
$ws.Range("T6").Value = 0.01297687277194152

$ws.Range("E7").Value = 3
$ws.Range("G7").Value = 152.3944216666667
$ws.Range("H7").Value = 457.183265
$ws.Range("I7").Value = 0.2369137948193439
$ws.Range("J7").Value = 0.2369137948193439
$ws.Range("K7").Value = 3
$ws.Range("M7").Value = 24.455837
$ws.Range("N7").Value = 73.36751100000001
$ws.Range("O7").Value = 0.1553502885444182
$ws.Range("P7").Value = 0.1553502885444182
$ws.Range("Q7").Value = 3726.933135989269
$ws.Range("R7").Value = 33542.39822390342
$ws.Range("S7").Value = 0.03680462638533816
$ws.Range("T7").Value = 0.03680462638533817

$ws.Range("E8").Value = 3
$ws.Range("G8").Value = 152.3944216666667
$ws.Range("H8").Value = 457.183265
$ws.Range("I8").Value = 0.2369137948193439
$ws.Range("J8").Value = 0.2369137948193439
$ws.Range("K8").Value = 3
$ws.Range("M8").Value = 33.36516466666667
$ws.Range("N8").Value = 100.095494
$ws.Range("O8").Value = 0.2119448194841458
$ws.Range("P8").Value = 0.2119448194841458
$ws.Range("Q8").Value = 5084.664973189767
$ws.Range("R8").Value = 45761.98475870791
$ws.Range("S8").Value = 0.05021265147628979
$ws.Range("T8").Value = 0.0502126514762898

$ws.Range("E9").Value = 3
$ws.Range("G9").Value = 152.3944216666667
$ws.Range("H9").Value = 457.183265
$ws.Range("I9").Value = 0.2369137948193439
$ws.Range("J9").Value = 0.2369137948193439
$ws.Range("K9").Value = 3
$ws.Range("M9").Value = 43.331795
$ws.Range("N9").Value = 129.995385
$ws.Range("O9").Value = 0.2752556314632608
$ws.Range("P9").Value = 0.2752556314632608
$ws.Range("Q9").Value = 6603.523838803558
$ws.Range("R9").Value = 59431.71454923203
$ws.Range("S9").Value = 0.06521185619535591
$ws.Range("T9").Value = 0.06521185619535591

$ws.Range("E10").Value = 3
$ws.Range("G10").Value = 152.3944216666667
$ws.Range("H10").Value = 457.183265
$ws.Range("I10").Value = 0.2369137948193439
$ws.Range("J10").Value = 0.2369137948193439
$ws.Range("K10").Value = 3
$ws.Range("M10").Value = 45.91844699999999
$ws.Range("N10").Value = 137.755341
$ws.Range("O10").Value = 0.2916867654524183
$ws.Range("P10").Value = 0.2916867654524183
$ws.Range("Q10").Value = 6997.715174396483
$ws.Range("R10").Value = 62979.43656956836
$ws.Range("S10").Value = 0.06910461850191231
$ws.Range("T10").Value = 0.06910461850191232

$ws.Range("E11").Value = 3
$ws.Range("G11").Value = 152.3944216666667
$ws.Range("H11").Value = 457.183265
$ws.Range("I11").Value = 0.2369137948193439
$ws.Range("J11").Value = 0.2369137948193439
$ws.Range("K11").Value = 3
$ws.Range("M11").Value = 10.35258366666667
$ws.Range("N11").Value = 31.057751
$ws.Range("O11").Value = 0.06576249505575693
$ws.Range("P11").Value = 0.06576249505575693
$ws.Range("Q11").Value = 1577.676000637446
$ws.Range("R11").Value = 14199.08400573702
$ws.Range("S11").Value = 0.01558004226044771
$ws.Range("T11").Value = 0.01558004226044772

$ws.Range("E12").Value = 3
$ws.Range("G12").Value = 206.573929
$ws.Range("H12").Value = 619.721787
$ws.Range("I12").Value = 0.321141764212203
$ws.Range("J12").Value = 0.321141764212203
$ws.Range("K12").Value = 3
$ws.Range("M12").Value = 24.455837
$ws.Range("N12").Value = 73.36751100000001
$ws.Range("O12").Value = 0.1553502885444182
$ws.Range("P12").Value = 0.1553502885444182
$ws.Range("Q12").Value = 5051.938336073574
$ws.Range("R12").Value = 45467.44502466216
$ws.Range("S12").Value = 0.04988946573402925
$ws.Range("T12").Value = 0.04988946573402925

$ws.Range("E13").Value = 3
$ws.Range("G13").Value = 206.573929
$ws.Range("H13").Value = 619.721787
$ws.Range("I13").Value = 0.321141764212203
$ws.Range("J13").Value = 0.321141764212203
$ws.Range("K13").Value = 3
$ws.Range("M13").Value = 33.36516466666667
$ws.Range("N13").Value = 100.095494
$ws.Range("O13").Value = 0.2119448194841458
$ws.Range("P13").Value = 0.2119448194841458
$ws.Range("Q13").Value = 6892.373156925308
$ws.Range("R13").Value = 62031.35841232778
$ws.Range("S13").Value = 0.06806433324477548
$ws.Range("T13").Value = 0.06806433324477548

$ws.Range("E14").Value = 3
$ws.Range("G14").Value = 206.573929
$ws.Range("H14").Value = 619.721787
$ws.Range("I14").Value = 0.321141764212203
$ws.Range("J14").Value = 0.321141764212203
$ws.Range("K14").Value = 3
$ws.Range("M14").Value = 43.331795
$ws.Range("N14").Value = 129.995385
$ws.Range("O14").Value = 0.2752556314632608
$ws.Range("P14").Value = 0.2752556314632608
$ws.Range("Q14").Value = 8951.219143772554
$ws.Range("R14").Value = 80560.97229395299
$ws.Range("S14").Value = 0.08839607909745556
$ws.Range("T14").Value = 0.08839607909745555

$ws.Range("E15").Value = 3
$ws.Range("G15").Value = 206.573929
$ws.Range("H15").Value = 619.721787
$ws.Range("I15").Value = 0.321141764212203
$ws.Range("J15").Value = 0.321141764212203
$ws.Range("K15").Value = 3
$ws.Range("M15").Value = 45.91844699999999
$ws.Range("N15").Value = 137.755341
$ws.Range("O15").Value = 0.2916867654524183
$ws.Range("P15").Value = 0.2916867654524183
$ws.Range("Q15").Value = 9485.554010368262
$ws.Range("R15").Value = 85369.98609331435
$ws.Range("S15").Value = 0.09367280245474068
$ws.Range("T15").Value = 0.09367280245474068

$ws.Range("E16").Value = 3
$ws.Range("G16").Value = 206.573929
$ws.Range("H16").Value = 619.721787
$ws.Range("I16").Value = 0.321141764212203
$ws.Range("J16").Value = 0.321141764212203
$ws.Range("K16").Value = 3
$ws.Range("M16").Value = 10.35258366666667
$ws.Range("N16").Value = 31.057751
$ws.Range("O16").Value = 0.06576249505575693
$ws.Range("P16").Value = 0.06576249505575693
$ws.Range("Q16").Value = 2138.57388332456
$ws.Range("R16").Value = 19247.16494992104
$ws.Range("S16").Value = 0.02111908368120206
$ws.Range("T16").Value = 0.02111908368120206

$ws.Range("E17").Value = 3
$ws.Range("G17").Value = 141.7744496666667
$ws.Range("H17").Value = 425.323349
$ws.Range("I17").Value = 0.2204038869114384
$ws.Range("J17").Value = 0.2204038869114385
$ws.Range("K17").Value = 3
$ws.Range("M17").Value = 24.455837
$ws.Range("N17").Value = 73.36751100000001
$ws.Range("O17").Value = 0.1553502885444182
$ws.Range("P17").Value = 0.1553502885444182
$ws.Range("Q17").Value = 3467.212831812705
$ws.Range("R17").Value = 31204.91548631434
$ws.Range("S17").Value = 0.03423980742800328
$ws.Range("T17").Value = 0.03423980742800329

$ws.Range("E18").Value = 3
$ws.Range("G18").Value = 141.7744496666667
$ws.Range("H18").Value = 425.323349
$ws.Range("I18").Value = 0.2204038869114384
$ws.Range("J18").Value = 0.2204038869114385
$ws.Range("K18").Value = 3
$ws.Range("M18").Value = 33.36516466666667
$ws.Range("N18").Value = 100.095494
$ws.Range("O18").Value = 0.2119448194841458
$ws.Range("P18").Value = 0.2119448194841458
$ws.Range("Q18").Value = 4730.327858654378
$ws.Range("R18").Value = 42572.95072788941
$ws.Range("S18").Value = 0.0467134620250489
$ws.Range("T18").Value = 0.04671346202504891

$ws.Range("E19").Value = 3
$ws.Range("G19").Value = 141.7744496666667
$ws.Range("H19").Value = 425.323349
$ws.Range("I19").Value = 0.2204038869114384
$ws.Range("J19").Value = 0.2204038869114385
$ws.Range("K19").Value = 3
$ws.Range("M19").Value = 43.331795
$ws.Range("N19").Value = 129.995385
$ws.Range("O19").Value = 0.2752556314632608
$ws.Range("P19").Value = 0.2752556314632608
$ws.Range("Q19").Value = 6143.341389193818
$ws.Range("R19").Value = 55290.07250274436
$ws.Range("S19").Value = 0.06066741106876512
$ws.Range("T19").Value = 0.06066741106876511

$ws.Range("E20").Value = 3
$ws.Range("G20").Value = 141.7744496666667
$ws.Range("H20").Value = 425.323349
$ws.Range("I20").Value = 0.2204038869114384
$ws.Range("J20").Value = 0.2204038869114385
$ws.Range("K20").Value = 3
$ws.Range("M20").Value = 45.91844699999999
$ws.Range("N20").Value = 137.755341
$ws.Range("O20").Value = 0.2916867654524183
$ws.Range("P20").Value = 0.2916867654524183
$ws.Range("Q20").Value = 6510.062552973
$ws.Range("R20").Value = 58590.56297675701
$ws.Range("S20").Value = 0.06428889686633807
$ws.Range("T20").Value = 0.06428889686633808

$ws.Range("E21").Value = 3
$ws.Range("G21").Value = 141.7744496666667
$ws.Range("H21").Value = 425.323349
$ws.Range("I21").Value = 0.2204038869114384
$ws.Range("J21").Value = 0.2204038869114385
$ws.Range("K21").Value = 3
$ws.Range("M21").Value = 10.35258366666667
$ws.Range("N21").Value = 31.057751
$ws.Range("O21").Value = 0.06576249505575693
$ws.Range("P21").Value = 0.06576249505575693
$ws.Range("Q21").Value = 1467.731851969789
$ws.Range("R21").Value = 13209.5866677281
$ws.Range("S21").Value = 0.01449430952328308
$ws.Range("T21").Value = 0.01449430952328308

$ws.Range("E22").Value = 3
$ws.Range("G22").Value = 15.57379533333333
$ws.Range("H22").Value = 46.721386
$ws.Range("I22").Value = 0.02421116804544314
$ws.Range("J22").Value = 0.02421116804544315
$ws.Range("K22").Value = 3
$ws.Range("M22").Value = 24.455837
$ws.Range("N22").Value = 73.36751100000001
$ws.Range("O22").Value = 0.1553502885444182
$ws.Range("P22").Value = 0.1553502885444182
$ws.Range("Q22").Value = 380.8702001433606
$ws.Range("R22").Value = 3427.831801290246
$ws.Range("S22").Value = 0.00376121194185699
$ws.Range("T22").Value = 0.003761211941856991

$ws.Range("E23").Value = 3
$ws.Range("G23").Value = 15.57379533333333
$ws.Range("H23").Value = 46.721386
$ws.Range("I23").Value = 0.02421116804544314
$ws.Range("J23").Value = 0.02421116804544315
$ws.Range("K23").Value = 3
$ws.Range("M23").Value = 33.36516466666667
$ws.Range("N23").Value = 100.095494
$ws.Range("O23").Value = 0.2119448194841458
$ws.Range("P23").Value = 0.2119448194841458
$ws.Range("Q23").Value = 519.6222457816315
$ws.Range("R23").Value = 4676.600212034684
$ws.Range("S23").Value = 0.005131431640891766
$ws.Range("T23").Value = 0.005131431640891766

$ws.Range("E24").Value = 3
$ws.Range("G24").Value = 15.57379533333333
$ws.Range("H24").Value = 46.721386
$ws.Range("I24").Value = 0.02421116804544314
$ws.Range("J24").Value = 0.02421116804544315
$ws.Range("K24").Value = 3
$ws.Range("M24").Value = 43.331795
$ws.Range("N24").Value = 129.995385
$ws.Range("O24").Value = 0.2752556314632608
$ws.Range("P24").Value = 0.2752556314632608
$ws.Range("Q24").Value = 674.8405067559565
$ws.Range("R24").Value = 6073.564560803609
$ws.Range("S24").Value = 0.006664260348811575
$ws.Range("T24").Value = 0.006664260348811574

$ws.Range("E25").Value = 3
$ws.Range("G25").Value = 15.57379533333333
$ws.Range("H25").Value = 46.721386
$ws.Range("I25").Value = 0.02421116804544314
$ws.Range("J25").Value = 0.02421116804544315
$ws.Range("K25").Value = 3
$ws.Range("M25").Value = 45.91844699999999
$ws.Range("N25").Value = 137.755341
$ws.Range("O25").Value = 0.2916867654524183
$ws.Range("P25").Value = 0.2916867654524183
$ws.Range("Q25").Value = 715.1244956025138
$ws.Range("R25").Value = 6436.120460422625
$ws.Range("S25").Value = 0.007062077295000259
$ws.Range("T25").Value = 0.007062077295000259

$ws.Range("E26").Value = 3
$ws.Range("G26").Value = 15.57379533333333
$ws.Range("H26").Value = 46.721386
$ws.Range("I26").Value = 0.02421116804544314
$ws.Range("J26").Value = 0.02421116804544315
$ws.Range("K26").Value = 3
$ws.Range("M26").Value = 10.35258366666667
$ws.Range("N26").Value = 31.057751
$ws.Range("O26").Value = 0.06576249505575693
$ws.Range("P26").Value = 0.06576249505575693
$ws.Range("Q26").Value = 161.2290191958762
$ws.Range("R26").Value = 1451.061172762886
$ws.Range("S26").Value = 0.001592186818882555
$ws.Range("T26").Value = 0.001592186818882555
